# Update metrics_19_8 results: retrain run shifted every model_19_8_N to
# model_19_8_(N+1 mod 25) and refreshed all computed metric columns (B:Q)
# with the new values produced by the retraining for the new LM (per
# commit "atualizado todo o treinamento para o novo lm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_19_8_24"
$ws.Cells.Item(2, 2).Value = 0.9992976945902621
$ws.Cells.Item(2, 3).Value = 0.9987169024240916
$ws.Cells.Item(2, 4).Value = 0.9990537023105441
$ws.Cells.Item(2, 5).Value = 0.9993361815137395
$ws.Cells.Item(2, 6).Value = 0.9992374210415716
$ws.Cells.Item(2, 7).Value = 0.0006555712723457703
$ws.Cells.Item(2, 8).Value = 0.001197715265636175
$ws.Cells.Item(2, 9).Value = 0.0008227493199864309
$ws.Cells.Item(2, 10).Value = 0.000941394591815568
$ws.Cells.Item(2, 11).Value = 0.0008820719559009994
$ws.Cells.Item(2, 12).Value = 0.002973796382833203
$ws.Cells.Item(2, 13).Value = 0.02560412608049277
$ws.Cells.Item(2, 14).Value = 1.001296563833362
$ws.Cells.Item(2, 15).Value = 0.02669414673558709
$ws.Cells.Item(2, 16).Value = 88.66000706185868
$ws.Cells.Item(2, 17).Value = 133.7584125819821

$ws.Cells.Item(3, 1).Value = "model_19_8_23"
$ws.Cells.Item(3, 2).Value = 0.9992835284735914
$ws.Cells.Item(3, 3).Value = 0.9987058758114503
$ws.Cells.Item(3, 4).Value = 0.9990340724140749
$ws.Cells.Item(3, 5).Value = 0.9993227908813958
$ws.Cells.Item(3, 6).Value = 0.9992218348551324
$ws.Cells.Item(3, 7).Value = 0.0006687947204372375
$ws.Cells.Item(3, 8).Value = 0.001208008124524552
$ws.Cells.Item(3, 9).Value = 0.0008398163425010184
$ws.Cells.Item(3, 10).Value = 0.0009603845252541853
$ws.Cells.Item(3, 11).Value = 0.0009001004338776018
$ws.Cells.Item(3, 12).Value = 0.003003740345307577
$ws.Cells.Item(3, 13).Value = 0.02586106572508638
$ws.Cells.Item(3, 14).Value = 1.001322716664139
$ws.Cells.Item(3, 15).Value = 0.02696202483278943
$ws.Cells.Item(3, 16).Value = 88.62006678061604
$ws.Cells.Item(3, 17).Value = 133.7184723007394

$ws.Cells.Item(4, 1).Value = "model_19_8_22"
$ws.Cells.Item(4, 2).Value = 0.9992690902744625
$ws.Cells.Item(4, 3).Value = 0.9986946800490661
$ws.Cells.Item(4, 4).Value = 0.9990140368275451
$ws.Cells.Item(4, 5).Value = 0.9993092124820709
$ws.Cells.Item(4, 6).Value = 0.9992059918353862
$ws.Cells.Item(4, 7).Value = 0.000682272145560302
$ws.Cells.Item(4, 8).Value = 0.001218458877273021
$ws.Cells.Item(4, 9).Value = 0.0008572360882920672
$ws.Cells.Item(4, 10).Value = 0.0009796407405515942
$ws.Cells.Item(4, 11).Value = 0.0009184259898879475
$ws.Cells.Item(4, 12).Value = 0.003033969200994667
$ws.Cells.Item(4, 13).Value = 0.02612033969075253
$ws.Cells.Item(4, 14).Value = 1.001349371800992
$ws.Cells.Item(4, 15).Value = 0.02723233662794509
$ws.Cells.Item(4, 16).Value = 88.58016387861052
$ws.Cells.Item(4, 17).Value = 133.6785693987339

$ws.Cells.Item(5, 1).Value = "model_19_8_21"
$ws.Cells.Item(5, 2).Value = 0.9992543836137502
$ws.Cells.Item(5, 3).Value = 0.9986832611933282
$ws.Cells.Item(5, 4).Value = 0.9989935796553536
$ws.Cells.Item(5, 5).Value = 0.9992954748044101
$ws.Cells.Item(5, 6).Value = 0.9991898712481494
$ws.Cells.Item(5, 7).Value = 0.0006960001677874212
$ws.Cells.Item(5, 8).Value = 0.001229117877874563
$ws.Cells.Item(5, 9).Value = 0.0008750223776351297
$ws.Cells.Item(5, 10).Value = 0.0009991228365187081
$ws.Cells.Item(5, 11).Value = 0.0009370726070769188
$ws.Cells.Item(5, 12).Value = 0.003064498390306193
$ws.Cells.Item(5, 13).Value = 0.02638181509652854
$ws.Cells.Item(5, 14).Value = 1.00137652255923
$ws.Cells.Item(5, 15).Value = 0.02750494358307369
$ws.Cells.Item(5, 16).Value = 88.540321313112
$ws.Cells.Item(5, 17).Value = 133.6387268332354

$ws.Cells.Item(6, 1).Value = "model_19_8_20"
$ws.Cells.Item(6, 2).Value = 0.999239393074694
$ws.Cells.Item(6, 3).Value = 0.9986715863865604
$ws.Cells.Item(6, 4).Value = 0.9989726788120102
$ws.Cells.Item(6, 5).Value = 0.9992814893865551
$ws.Cells.Item(6, 6).Value = 0.999173446971876
$ws.Cells.Item(6, 7).Value = 0.0007099931779877862
$ws.Cells.Item(6, 8).Value = 0.001240015797527502
$ws.Cells.Item(6, 9).Value = 0.0008931944125449394
$ws.Cells.Item(6, 10).Value = 0.001018956265393312
$ws.Cells.Item(6, 11).Value = 0.000956070500129865
$ws.Cells.Item(6, 12).Value = 0.003095303354256647
$ws.Cells.Item(6, 13).Value = 0.02664569717586287
$ws.Cells.Item(6, 14).Value = 1.001404197400565
$ws.Cells.Item(6, 15).Value = 0.02778005966883649
$ws.Cells.Item(6, 16).Value = 88.50051039288596
$ws.Cells.Item(6, 17).Value = 133.5989159130094

$ws.Cells.Item(7, 1).Value = "model_19_8_19"
$ws.Cells.Item(7, 2).Value = 0.9992241236384317
$ws.Cells.Item(7, 3).Value = 0.9986596730435853
$ws.Cells.Item(7, 4).Value = 0.9989513303462617
$ws.Cells.Item(7, 5).Value = 0.9992672891572829
$ws.Cells.Item(7, 6).Value = 0.9991567144161128
$ws.Cells.Item(7, 7).Value = 0.0007242465264878399
$ws.Cells.Item(7, 8).Value = 0.001251136380259432
$ws.Cells.Item(7, 9).Value = 0.0009117556284002285
$ws.Cells.Item(7, 10).Value = 0.001039094329210679
$ws.Cells.Item(7, 11).Value = 0.0009754249788054539
$ws.Cells.Item(7, 12).Value = 0.003126386883239584
$ws.Cells.Item(7, 13).Value = 0.02691182874662812
$ws.Cells.Item(7, 14).Value = 1.001432387129049
$ws.Cells.Item(7, 15).Value = 0.02805752101153752
$ws.Cells.Item(7, 16).Value = 88.46075743467793
$ws.Cells.Item(7, 17).Value = 133.5591629548013

$ws.Cells.Item(8, 1).Value = "model_19_8_18"
$ws.Cells.Item(8, 2).Value = 0.999208572361741
$ws.Cells.Item(8, 3).Value = 0.9986475447261876
$ws.Cells.Item(8, 4).Value = 0.9989295357588168
$ws.Cells.Item(8, 5).Value = 0.9992529211554511
$ws.Cells.Item(8, 6).Value = 0.9991397224528374
$ws.Cells.Item(8, 7).Value = 0.000738762960656539
$ws.Cells.Item(8, 8).Value = 0.001262457632178593
$ws.Cells.Item(8, 9).Value = 0.000930704720424351
$ws.Cells.Item(8, 10).Value = 0.0010594703197859
$ws.Cells.Item(8, 11).Value = 0.0009950795130870826
$ws.Cells.Item(8, 12).Value = 0.003157738921295777
$ws.Cells.Item(8, 13).Value = 0.02718019427186898
$ws.Cells.Item(8, 14).Value = 1.001461097178324
$ws.Cells.Item(8, 15).Value = 0.0283373114127811
$ws.Cells.Item(8, 16).Value = 88.42106689072388
$ws.Cells.Item(8, 17).Value = 133.5194724108473

$ws.Cells.Item(9, 1).Value = "model_19_8_17"
$ws.Cells.Item(9, 2).Value = 0.9991927413071447
$ws.Cells.Item(9, 3).Value = 0.9986351480617069
$ws.Cells.Item(9, 4).Value = 0.9989073058171645
$ws.Cells.Item(9, 5).Value = 0.9992383041316842
$ws.Cells.Item(9, 6).Value = 0.9991224003711789
$ws.Cells.Item(9, 7).Value = 0.0007535405552190452
$ws.Cells.Item(9, 8).Value = 0.001274029374320562
$ws.Cells.Item(9, 9).Value = 0.0009500323269287511
$ws.Cells.Item(9, 10).Value = 0.001080199460970343
$ws.Cells.Item(9, 11).Value = 0.001015115893949547
$ws.Cells.Item(9, 12).Value = 0.00318937167886259
$ws.Cells.Item(9, 13).Value = 0.02745069316463694
$ws.Cells.Item(9, 14).Value = 1.001490323740656
$ws.Cells.Item(9, 15).Value = 0.02861932600342409
$ws.Cells.Item(9, 16).Value = 88.38145543778431
$ws.Cells.Item(9, 17).Value = 133.4798609579077

$ws.Cells.Item(10, 1).Value = "model_19_8_16"
$ws.Cells.Item(10, 2).Value = 0.9991766122746095
$ws.Cells.Item(10, 3).Value = 0.9986225524130879
$ws.Cells.Item(10, 4).Value = 0.9988846100093703
$ws.Cells.Item(10, 5).Value = 0.9992235189520785
$ws.Cells.Item(10, 6).Value = 0.9991048070431828
$ws.Cells.Item(10, 7).Value = 0.0007685962991080155
$ws.Cells.Item(10, 8).Value = 0.00128578685942134
$ws.Cells.Item(10, 9).Value = 0.0009697649762178619
$ws.Cells.Item(10, 10).Value = 0.001101167072460485
$ws.Cells.Item(10, 11).Value = 0.001035466024339174
$ws.Cells.Item(10, 12).Value = 0.003221292496966393
$ws.Cells.Item(10, 13).Value = 0.02772356937892406
$ws.Cells.Item(10, 14).Value = 1.001520100416105
$ws.Cells.Item(10, 15).Value = 0.02890381912308505
$ws.Cells.Item(10, 16).Value = 88.34188938990619
$ws.Cells.Item(10, 17).Value = 133.4402949100296

$ws.Cells.Item(11, 1).Value = "model_19_8_15"
$ws.Cells.Item(11, 2).Value = 0.999160168920147
$ws.Cells.Item(11, 3).Value = 0.9986097006979058
$ws.Cells.Item(11, 4).Value = 0.998861383109357
$ws.Cells.Item(11, 5).Value = 0.9992085009228137
$ws.Cells.Item(11, 6).Value = 0.9990868669655927
$ws.Cells.Item(11, 7).Value = 0.0007839454487189648
$ws.Cells.Item(11, 8).Value = 0.001297783371418746
$ws.Cells.Item(11, 9).Value = 0.0009899593784702099
$ws.Cells.Item(11, 10).Value = 0.00112246489983683
$ws.Cells.Item(11, 11).Value = 0.001056217238563013
$ws.Cells.Item(11, 12).Value = 0.003253528611249728
$ws.Cells.Item(11, 13).Value = 0.02799902585303576
$ws.Cells.Item(11, 14).Value = 1.00155045737819
$ws.Cells.Item(11, 15).Value = 0.0291910023495732
$ws.Cells.Item(11, 16).Value = 88.30234224150064
$ws.Cells.Item(11, 17).Value = 133.4007477616241

$ws.Cells.Item(12, 1).Value = "model_19_8_14"
$ws.Cells.Item(12, 2).Value = 0.9991434430793603
$ws.Cells.Item(12, 3).Value = 0.9985966422543754
$ws.Cells.Item(12, 4).Value = 0.9988376456029966
$ws.Cells.Item(12, 5).Value = 0.9991933640507314
$ws.Cells.Item(12, 6).Value = 0.9990686754058385
$ws.Cells.Item(12, 7).Value = 0.0007995582869137915
$ws.Cells.Item(12, 8).Value = 0.001309972855254943
$ws.Cells.Item(12, 9).Value = 0.001010597722443653
$ws.Cells.Item(12, 10).Value = 0.001143931264227326
$ws.Cells.Item(12, 11).Value = 0.001077259341175304
$ws.Cells.Item(12, 12).Value = 0.003286011682170793
$ws.Cells.Item(12, 13).Value = 0.02827646171135617
$ws.Cells.Item(12, 14).Value = 1.001581335853489
$ws.Cells.Item(12, 15).Value = 0.02948024922675371
$ws.Cells.Item(12, 16).Value = 88.26290224828055
$ws.Cells.Item(12, 17).Value = 133.361307768404

$ws.Cells.Item(13, 1).Value = "model_19_8_13"
$ws.Cells.Item(13, 2).Value = 0.9991263890413723
$ws.Cells.Item(13, 3).Value = 0.9985832740329849
$ws.Cells.Item(13, 4).Value = 0.9988134228269502
$ws.Cells.Item(13, 5).Value = 0.9991779523803291
$ws.Cells.Item(13, 6).Value = 0.9990501196942484
$ws.Cells.Item(13, 7).Value = 0.0008154774827898344
$ws.Cells.Item(13, 8).Value = 0.001322451503125898
$ws.Cells.Item(13, 9).Value = 0.001031657979424515
$ws.Cells.Item(13, 10).Value = 0.001165787334023277
$ws.Cells.Item(13, 11).Value = 0.001098722656723896
$ws.Cells.Item(13, 12).Value = 0.003318823018824734
$ws.Cells.Item(13, 13).Value = 0.02855656636904784
$ws.Cells.Item(13, 14).Value = 1.001612820231313
$ws.Cells.Item(13, 15).Value = 0.02977227851961989
$ws.Cells.Item(13, 16).Value = 88.2234734956471
$ws.Cells.Item(13, 17).Value = 133.3218790157705

$ws.Cells.Item(14, 1).Value = "model_19_8_12"
$ws.Cells.Item(14, 2).Value = 0.9991090009827165
$ws.Cells.Item(14, 3).Value = 0.9985696199525095
$ws.Cells.Item(14, 4).Value = 0.9987886954607899
$ws.Cells.Item(14, 5).Value = 0.9991622824572169
$ws.Cells.Item(14, 6).Value = 0.9990312204836095
$ws.Cells.Item(14, 7).Value = 0.00083170847229745
$ws.Cells.Item(14, 8).Value = 0.001335196987904813
$ws.Cells.Item(14, 9).Value = 0.00105315694737101
$ws.Cells.Item(14, 10).Value = 0.001188009645057577
$ws.Cells.Item(14, 11).Value = 0.001120583296214294
$ws.Cells.Item(14, 12).Value = 0.003351973960686426
$ws.Cells.Item(14, 13).Value = 0.02883935630865311
$ws.Cells.Item(14, 14).Value = 1.001644921262677
$ws.Cells.Item(14, 15).Value = 0.03006710741241003
$ws.Cells.Item(14, 16).Value = 88.1840571448372
$ws.Cells.Item(14, 17).Value = 133.2824626649606

$ws.Cells.Item(15, 1).Value = "model_19_8_11"
$ws.Cells.Item(15, 2).Value = 0.9990912813626237
$ws.Cells.Item(15, 3).Value = 0.9985556861145609
$ws.Cells.Item(15, 4).Value = 0.9987634251137577
$ws.Cells.Item(15, 5).Value = 0.9991463955119151
$ws.Cells.Item(15, 6).Value = 0.9990119841663957
$ws.Cells.Item(15, 7).Value = 0.0008482489598526698
$ws.Cells.Item(15, 8).Value = 0.001348203614004987
$ws.Cells.Item(15, 9).Value = 0.001075128004754159
$ws.Cells.Item(15, 10).Value = 0.00121053972624264
$ws.Cells.Item(15, 11).Value = 0.001142833865498399
$ws.Cells.Item(15, 12).Value = 0.003385436141524916
$ws.Cells.Item(15, 13).Value = 0.02912471390164494
$ws.Cells.Item(15, 14).Value = 1.001677634407464
$ws.Cells.Item(15, 15).Value = 0.03036461326890717
$ws.Cells.Item(15, 16).Value = 88.14467276105603
$ws.Cells.Item(15, 17).Value = 133.2430782811794

$ws.Cells.Item(16, 1).Value = "model_19_8_10"
$ws.Cells.Item(16, 2).Value = 0.9990732382973089
$ws.Cells.Item(16, 3).Value = 0.9985414975940605
$ws.Cells.Item(16, 4).Value = 0.9987375663577863
$ws.Cells.Item(16, 5).Value = 0.9991303379561075
$ws.Cells.Item(16, 6).Value = 0.9989924221206047
$ws.Cells.Item(16, 7).Value = 0.0008650913693250295
$ws.Cells.Item(16, 8).Value = 0.001361447975088075
$ws.Cells.Item(16, 9).Value = 0.001097610648565027
$ws.Cells.Item(16, 10).Value = 0.001233311758820692
$ws.Cells.Item(16, 11).Value = 0.001165461203692859
$ws.Cells.Item(16, 12).Value = 0.003419189818530003
$ws.Cells.Item(16, 13).Value = 0.02941243562381445
$ws.Cells.Item(16, 14).Value = 1.001710944681891
$ws.Cells.Item(16, 15).Value = 0.03066458390045551
$ws.Cells.Item(16, 16).Value = 88.10535085466657
$ws.Cells.Item(16, 17).Value = 133.20375637479

$ws.Cells.Item(17, 1).Value = "model_19_8_9"
$ws.Cells.Item(17, 2).Value = 0.9990548362304459
$ws.Cells.Item(17, 3).Value = 0.9985270150299962
$ws.Cells.Item(17, 4).Value = 0.9987110767133563
$ws.Cells.Item(17, 5).Value = 0.9991140226362842
$ws.Cells.Item(17, 6).Value = 0.9989724649545284
$ws.Cells.Item(17, 7).Value = 0.0008822688909842776
$ws.Cells.Item(17, 8).Value = 0.001374966812931007
$ws.Cells.Item(17, 9).Value = 0.001120641812208693
$ws.Cells.Item(17, 10).Value = 0.001256449339595158
$ws.Cells.Item(17, 11).Value = 0.001188545575901925
$ws.Cells.Item(17, 12).Value = 0.003453285873252217
$ws.Cells.Item(17, 13).Value = 0.02970301148005498
$ws.Cells.Item(17, 14).Value = 1.001744917728408
$ws.Cells.Item(17, 15).Value = 0.03096753017247114
$ws.Cells.Item(17, 16).Value = 88.0660273667157
$ws.Cells.Item(17, 17).Value = 133.1644328868391

$ws.Cells.Item(18, 1).Value = "model_19_8_8"
$ws.Cells.Item(18, 2).Value = 0.9990361088404808
$ws.Cells.Item(18, 3).Value = 0.998512226886777
$ws.Cells.Item(18, 4).Value = 0.9986841006126645
$ws.Cells.Item(18, 5).Value = 0.9990974705546007
$ws.Cells.Item(18, 6).Value = 0.9989521798247064
$ws.Cells.Item(18, 7).Value = 0.0008997500874792026
$ws.Cells.Item(18, 8).Value = 0.001388770895501722
$ws.Cells.Item(18, 9).Value = 0.001144095920516658
$ws.Cells.Item(18, 10).Value = 0.001279922684346238
$ws.Cells.Item(18, 11).Value = 0.001212009302431448
$ws.Cells.Item(18, 12).Value = 0.003487674644430308
$ws.Cells.Item(18, 13).Value = 0.0299958345021305
$ws.Cells.Item(18, 14).Value = 1.00177949137142
$ws.Cells.Item(18, 15).Value = 0.03127281927682364
$ws.Cells.Item(18, 16).Value = 88.02678702755803
$ws.Cells.Item(18, 17).Value = 133.1251925476815

$ws.Cells.Item(19, 1).Value = "model_19_8_7"
$ws.Cells.Item(19, 2).Value = 0.9990169907552966
$ws.Cells.Item(19, 3).Value = 0.9984970961363228
$ws.Cells.Item(19, 4).Value = 0.9986564489166853
$ws.Cells.Item(19, 5).Value = 0.9990806136558952
$ws.Cells.Item(19, 6).Value = 0.9989314539281304
$ws.Cells.Item(19, 7).Value = 0.0009175959808117312
$ws.Cells.Item(19, 8).Value = 0.001402894786887567
$ws.Cells.Item(19, 9).Value = 0.001168137418574695
$ws.Cells.Item(19, 10).Value = 0.001303828305543322
$ws.Cells.Item(19, 11).Value = 0.001235982862059009
$ws.Cells.Item(19, 12).Value = 0.003522438372665244
$ws.Cells.Item(19, 13).Value = 0.03029184677123089
$ws.Cells.Item(19, 14).Value = 1.001814786297914
$ws.Cells.Item(19, 15).Value = 0.03158143340104951
$ws.Cells.Item(19, 16).Value = 87.98750674453596
$ws.Cells.Item(19, 17).Value = 133.0859122646594

$ws.Cells.Item(20, 1).Value = "model_19_8_6"
$ws.Cells.Item(20, 2).Value = 0.9989975158999641
$ws.Cells.Item(20, 3).Value = 0.9984816548329628
$ws.Cells.Item(20, 4).Value = 0.9986282349623683
$ws.Cells.Item(20, 5).Value = 0.9990635448396622
$ws.Cells.Item(20, 6).Value = 0.9989103868091935
$ws.Cells.Item(20, 7).Value = 0.0009357749034172982
$ws.Cells.Item(20, 8).Value = 0.001417308565779186
$ws.Cells.Item(20, 9).Value = 0.001192667766674539
$ws.Cells.Item(20, 10).Value = 0.001328034457711417
$ws.Cells.Item(20, 11).Value = 0.001260351112192978
$ws.Cells.Item(20, 12).Value = 0.003557547136315668
$ws.Cells.Item(20, 13).Value = 0.0305904381043701
$ws.Cells.Item(20, 14).Value = 1.001850739876989
$ws.Cells.Item(20, 15).Value = 0.03189273638540972
$ws.Cells.Item(20, 16).Value = 87.94827119642109
$ws.Cells.Item(20, 17).Value = 133.0466767165445

$ws.Cells.Item(21, 1).Value = "model_19_8_5"
$ws.Cells.Item(21, 2).Value = 0.9989776612226549
$ws.Cells.Item(21, 3).Value = 0.9984658978691939
$ws.Cells.Item(21, 4).Value = 0.9985993653254582
$ws.Cells.Item(21, 5).Value = 0.9990461774048406
$ws.Cells.Item(21, 6).Value = 0.9988888902070259
$ws.Cells.Item(21, 7).Value = 0.0009543083731658068
$ws.Cells.Item(21, 8).Value = 0.001432016999806606
$ws.Cells.Item(21, 9).Value = 0.001217768191626124
$ws.Cells.Item(21, 10).Value = 0.001352664096013408
$ws.Cells.Item(21, 11).Value = 0.001285216143819766
$ws.Cells.Item(21, 12).Value = 0.003592986212007804
$ws.Cells.Item(21, 13).Value = 0.03089188199455978
$ws.Cells.Item(21, 14).Value = 1.001887394665868
$ws.Cells.Item(21, 15).Value = 0.03220701336607966
$ws.Cells.Item(21, 16).Value = 87.90904739287072
$ws.Cells.Item(21, 17).Value = 133.0074529129942

$ws.Cells.Item(22, 1).Value = "model_19_8_4"
$ws.Cells.Item(22, 2).Value = 0.9989573747727022
$ws.Cells.Item(22, 3).Value = 0.9984497412527488
$ws.Cells.Item(22, 4).Value = 0.9985697500299795
$ws.Cells.Item(22, 5).Value = 0.9990284741751165
$ws.Cells.Item(22, 6).Value = 0.9988669075162793
$ws.Cells.Item(22, 7).Value = 0.0009732448837243585
$ws.Cells.Item(22, 8).Value = 0.001447098492064658
$ws.Cells.Item(22, 9).Value = 0.001243516922165976
$ws.Cells.Item(22, 10).Value = 0.001377769942061376
$ws.Cells.Item(22, 11).Value = 0.001310643432113676
$ws.Cells.Item(22, 12).Value = 0.003628859059660633
$ws.Cells.Item(22, 13).Value = 0.03119687297990551
$ws.Cells.Item(22, 14).Value = 1.001924846573473
$ws.Cells.Item(22, 15).Value = 0.03252498844908993
$ws.Cells.Item(22, 16).Value = 87.86974965677037
$ws.Cells.Item(22, 17).Value = 132.9681551768938

$ws.Cells.Item(23, 1).Value = "model_19_8_3"
$ws.Cells.Item(23, 2).Value = 0.9989367473939437
$ws.Cells.Item(23, 3).Value = 0.9984332773404617
$ws.Cells.Item(23, 4).Value = 0.9985395195587231
$ws.Cells.Item(23, 5).Value = 0.9990105620897937
$ws.Cells.Item(23, 6).Value = 0.9988445655922678
$ws.Cells.Item(23, 7).Value = 0.000992499636358228
$ws.Cells.Item(23, 8).Value = 0.001462466831502415
$ws.Cells.Item(23, 9).Value = 0.001269800511301
$ws.Cells.Item(23, 10).Value = 0.001403171976804587
$ws.Cells.Item(23, 11).Value = 0.001336486244052794
$ws.Cells.Item(23, 12).Value = 0.003664990777370118
$ws.Cells.Item(23, 13).Value = 0.03150396223268159
$ws.Cells.Item(23, 14).Value = 1.001962927888104
$ws.Cells.Item(23, 15).Value = 0.03284515112711907
$ws.Cells.Item(23, 16).Value = 87.83056782358538
$ws.Cells.Item(23, 17).Value = 132.9289733437088

$ws.Cells.Item(24, 1).Value = "model_19_8_2"
$ws.Cells.Item(24, 2).Value = 0.9989156801594219
$ws.Cells.Item(24, 3).Value = 0.9984164263240463
$ws.Cells.Item(24, 4).Value = 0.9985085767762117
$ws.Cells.Item(24, 5).Value = 0.9989923029377588
$ws.Cells.Item(24, 6).Value = 0.9988217432027516
$ws.Cells.Item(24, 7).Value = 0.00101216497503955
$ws.Cells.Item(24, 8).Value = 0.001478196515652165
$ws.Cells.Item(24, 9).Value = 0.001296703412526906
$ws.Cells.Item(24, 10).Value = 0.001429066204417321
$ws.Cells.Item(24, 11).Value = 0.001362884808472113
$ws.Cells.Item(24, 12).Value = 0.00370156538864877
$ws.Cells.Item(24, 13).Value = 0.03181454030847453
$ws.Cells.Item(24, 14).Value = 1.002001821244144
$ws.Cells.Item(24, 15).Value = 0.03316895115458376
$ws.Cells.Item(24, 16).Value = 87.79132740517841
$ws.Cells.Item(24, 17).Value = 132.8897329253018

$ws.Cells.Item(25, 1).Value = "model_19_8_1"
$ws.Cells.Item(25, 2).Value = 0.9988941771112638
$ws.Cells.Item(25, 3).Value = 0.9983991728297686
$ws.Cells.Item(25, 4).Value = 0.9984769177518696
$ws.Cells.Item(25, 5).Value = 0.9989737150965917
$ws.Cells.Item(25, 6).Value = 0.9987984450554158
$ws.Cells.Item(25, 7).Value = 0.001032237126620414
$ws.Cells.Item(25, 8).Value = 0.001494301895219427
$ws.Cells.Item(25, 9).Value = 0.001324229043244456
$ws.Cells.Item(25, 10).Value = 0.001455426562723627
$ws.Cells.Item(25, 11).Value = 0.001389833680011505
$ws.Cells.Item(25, 12).Value = 0.003738543111039413
$ws.Cells.Item(25, 13).Value = 0.03212844731107332
$ws.Cells.Item(25, 14).Value = 1.002041519179205
$ws.Cells.Item(25, 15).Value = 0.03349622182816026
$ws.Cells.Item(25, 16).Value = 87.75205372891949
$ws.Cells.Item(25, 17).Value = 132.8504592490429

$ws.Cells.Item(26, 1).Value = "model_19_8_0"
$ws.Cells.Item(26, 2).Value = 0.9988723031439045
$ws.Cells.Item(26, 3).Value = 0.9983816091630108
$ws.Cells.Item(26, 4).Value = 0.9984445610516192
$ws.Cells.Item(26, 5).Value = 0.9989549013102997
$ws.Cells.Item(26, 6).Value = 0.9987747563544189
$ws.Cells.Item(26, 7).Value = 0.001052655514994088
$ws.Cells.Item(26, 8).Value = 0.001510696807182013
$ws.Cells.Item(26, 9).Value = 0.001352361261493305
$ws.Cells.Item(26, 10).Value = 0.001482107345246822
$ws.Cells.Item(26, 11).Value = 0.001417234303370064
$ws.Cells.Item(26, 12).Value = 0.003775792151771995
$ws.Cells.Item(26, 13).Value = 0.03244465310330946
$ws.Cells.Item(26, 14).Value = 1.002081901888176
$ws.Cells.Item(26, 15).Value = 0.03382588915560811
$ws.Cells.Item(26, 16).Value = 87.71287849121747
$ws.Cells.Item(26, 17).Value = 132.8112840113409
